$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at position 12 (pushes the two "Order" rows down by one,
# i.e. old row 13 -> 14, old row 14 -> 15; the previously-blank row 12 -> 13)
$ws.Rows.Item(12).Insert()

# New row 12: a "View" rule for the Free Cash Flow to Firm Growth (fcffgrowth) metric
$ws.Range("A12").Value = "View"
$ws.Range("B12").Value = "Standard"
$ws.Range("C12").Value = "Free Cash Flow to Firm"
$ws.Range("D12").Value = "fcffgrowth"
$ws.Range("E12").Value = "~gt~"
$ws.Range("F12").Value = -99

# Row 13 becomes the former "Order/Risky/Revenue Growth/desc" row, now expressed
# as a "View" rule instead of an "Order" rule
$ws.Range("A13").Value = "View"
$ws.Range("B13").Value = "Risky"
$ws.Range("C13").Value = "Revenue Growth"
$ws.Range("D13").Value = "revenuegrowth"
$ws.Range("E13").Value = "desc"

# The old data that used to live in rows 13/14 (now shifted to 14/15) is no
# longer needed since it has been folded into rows 12/13 above - remove it
$ws.Range("A14:A15").EntireRow.Delete()

# Match the saved selection state from the edit
$ws.Range("B13").Select() | Out-Null
